# Daily update of covid19 tracker data files
# Bump the "Date Added" column (B) for every data row on the "Country
# Updates" sheet from 2020-04-13 (serial 43934) to 2020-04-14 (serial 43935).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country Updates")

for ($row = 5; $row -le 96; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    if ($cell.Value2 -eq 43934) {
        $cell.Value2 = 43935
    }
}
